$wb = $excel.ActiveWorkbook

# --- Sheet1 (Analysis_Unit) / Sheet2 (r AnalysisUnit_Variable) handles ---
$ws1 = $wb.Worksheets.Item("Analysis_Unit")
$ws2 = $wb.Worksheets.Item("r AnalysisUnit_Variable")

$suffixes = @(52,53,54,55,60,61,62,63,64,70,71,80,81,82,83,84,85,86,87,100,114,201,202,203,204,205,206,207,208,209,210,211,212,213,214,215,216,217,218,219,220,221,222,223,224,225,226,227,228,229,230,231,106,107,108,109,110,111,112,113,115,116,117)

$startRow = 26
# First fill column B and C for every new row (this is the order the shared
# strings were originally added in: all COUNTERPARTY_RETAIL_IND_* first)
for ($i = 0; $i -lt $suffixes.Length; $i++) {
    $row = $startRow + $i
    $suffix = $suffixes[$i]
    if ($suffix -lt 100) {
        $cWord = "COUNTERPARTY_RETAIL_IND_$suffix "
    } else {
        $cWord = "COUNTERPARTY_RETAIL_IND_$suffix"
    }
    $ws2.Cells.Item($row, 2).Value = $cWord
    $ws2.Cells.Item($row, 3).Value = $cWord
}

# Then fill the rest of the row (A, E, F) for every new row, F getting the
# RETAIL_IND_* (non-counterparty) variant, in the same order
for ($i = 0; $i -lt $suffixes.Length; $i++) {
    $row = $startRow + $i
    $suffix = $suffixes[$i]
    if ($suffix -lt 100) {
        $fWord = "RETAIL_IND_$suffix "
    } else {
        $fWord = "RETAIL_IND_$suffix"
    }
    $ws2.Cells.Item($row, 1).Value = "CREATE/MODIFY"
    $ws2.Cells.Item($row, 5).Value = "COUNTERPARTY_RETAIL"
    $ws2.Cells.Item($row, 6).Value = $fWord
}


# --- Restore selections on both sheets; leave "Analysis_Unit" as the active
# (tabSelected) sheet, matching the original file's tab state ---
$ws2.Activate()
$ws2.Range("B96").Select()

$ws1.Activate()
$ws1.Range("B16").Select()
